{"js": "// Update the worksheet date and all 25 \"NNN\u00d7N=\" multiplication prompts to\n// the new values given in the commit's diff. Every <w:t> run in the document\n// is a unique string, so a plain exact-text search+replace (in document\n// order) reproduces the diff without touching unrelated whitespace/runs.\nconst replacements = [\n  [\"2024-04-22 Monday\", \"2024-04-23 Tuesday\"],\n  [\"982\u00d77=\", \"922\u00d76=\"],\n  [\"136\u00d77=\", \"584\u00d73=\"],\n  [\"388\u00d78=\", \"964\u00d79=\"],\n  [\"156\u00d73=\", \"564\u00d79=\"],\n  [\"601\u00d72=\", \"683\u00d72=\"],\n  [\"262\u00d72=\", \"418\u00d79=\"],\n  [\"536\u00d73=\", \"617\u00d77=\"],\n  [\"870\u00d77=\", \"846\u00d72=\"],\n  [\"250\u00d76=\", \"401\u00d78=\"],\n  [\"267\u00d78=\", \"623\u00d79=\"],\n  [\"361\u00d77=\", \"924\u00d74=\"],\n  [\"877\u00d72=\", \"406\u00d77=\"],\n  [\"201\u00d76=\", \"123\u00d72=\"],\n  [\"741\u00d75=\", \"931\u00d74=\"],\n  [\"222\u00d79=\", \"211\u00d77=\"],\n  [\"885\u00d77=\", \"737\u00d76=\"],\n  [\"298\u00d74=\", \"536\u00d79=\"],\n  [\"539\u00d79=\", \"271\u00d72=\"],\n  [\"428\u00d73=\", \"353\u00d74=\"],\n  [\"429\u00d78=\", \"356\u00d77=\"],\n  [\"315\u00d76=\", \"855\u00d75=\"],\n  [\"612\u00d74=\", \"817\u00d73=\"],\n  [\"309\u00d76=\", \"305\u00d73=\"],\n  [\"337\u00d72=\", \"184\u00d74=\"],\n  [\"820\u00d74=\", \"911\u00d78=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: ${oldText}`);\n  }\n\n  for (const item of results.items) {\n    item.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Update the worksheet date and all 25 \"NNN\u00d7N=\" multiplication prompts to\n# the new values given in the commit's diff. Every run in the document holds\n# a unique string, so a plain exact-text Find/Replace (one pair at a time, in\n# document order) reproduces the diff without touching unrelated\n# whitespace/runs or formatting.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2024-04-22 Monday\", \"2024-04-23 Tuesday\"),\n    @(\"982\u00d77=\", \"922\u00d76=\"),\n    @(\"136\u00d77=\", \"584\u00d73=\"),\n    @(\"388\u00d78=\", \"964\u00d79=\"),\n    @(\"156\u00d73=\", \"564\u00d79=\"),\n    @(\"601\u00d72=\", \"683\u00d72=\"),\n    @(\"262\u00d72=\", \"418\u00d79=\"),\n    @(\"536\u00d73=\", \"617\u00d77=\"),\n    @(\"870\u00d77=\", \"846\u00d72=\"),\n    @(\"250\u00d76=\", \"401\u00d78=\"),\n    @(\"267\u00d78=\", \"623\u00d79=\"),\n    @(\"361\u00d77=\", \"924\u00d74=\"),\n    @(\"877\u00d72=\", \"406\u00d77=\"),\n    @(\"201\u00d76=\", \"123\u00d72=\"),\n    @(\"741\u00d75=\", \"931\u00d74=\"),\n    @(\"222\u00d79=\", \"211\u00d77=\"),\n    @(\"885\u00d77=\", \"737\u00d76=\"),\n    @(\"298\u00d74=\", \"536\u00d79=\"),\n    @(\"539\u00d79=\", \"271\u00d72=\"),\n    @(\"428\u00d73=\", \"353\u00d74=\"),\n    @(\"429\u00d78=\", \"356\u00d77=\"),\n    @(\"315\u00d76=\", \"855\u00d75=\"),\n    @(\"612\u00d74=\", \"817\u00d73=\"),\n    @(\"309\u00d76=\", \"305\u00d73=\"),\n    @(\"337\u00d72=\", \"184\u00d74=\"),\n    @(\"820\u00d74=\", \"911\u00d78=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Forward = $true\n    $find.Wrap = 1\n    $find.Format = $false\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n\n    $found = $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n\n    if (-not $found) {\n        throw \"Text not found: $oldText\"\n    }\n}\n"}
